# Update to perk database
# Row 27: "Cover" Repair Bolt Speed (repair_bolt) bonuses change from
#   10%-PROJECTILE_SPEED / 25%-INVULNERABILITY  ->  5%-PROJECTILE_SPEED / 5%-INVULNERABILITY
# Row 33: "Speedy" Repairing Cooldown (repair_bolt) 2nd bonus changes from
#   13%-PROJECTILE_SPEED -> 5%-PROJECTILE_SPEED
# Both rows also pick up the slightly shaded formatting already used on row 31
# (style index 13/8/9 instead of 12/7/7), and the sheet's active selection
# moves to B24 with the view scrolled down a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 ------------------------------------------------------------
$ws.Range("C27").Value = "5%-PROJECTILE_SPEED"
$ws.Range("D27").Value = "5%-INVULNERABILITY"

# Re-stripe the row to match the shaded style already used elsewhere
# (e.g. row 31) by copying that row's cell formatting across.
$ws.Range("A31").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("B31").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("D31").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D31").Copy()
$ws.Range("D27").PasteSpecial(-4122)

# --- Row 33 --------------------------------------------------------------
$ws.Range("D33").Value = "5%-PROJECTILE_SPEED"

$ws.Range("A31").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("B31").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("D28").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("D31").Copy()
$ws.Range("D33").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Selection / view ----------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("B24").Select() | Out-Null
